# Update timelines for CHUT
# - Fix the optional_lookback value for the CHUT DAP (was mistakenly entered
#   as 75625 instead of 75.625, matching the style of the other studies'
#   lookback/coverage values).
# - Apply the "Text" (@ ) number format consistently across columns
#   A, B, G, H and I (mirrors the formatting already used on columns C-F).
# - Keep the highlight fill on column B (DAP codes) combined with the new
#   text format.
# - Move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the CHUT optional_lookback values (column G) -------------------
# Must happen BEFORE the column is switched to Text format, otherwise the
# numeric value would be captured as a text string instead of a number.
$ws.Range("G2").Value = 75.625
$ws.Range("G8").Value = 75.625
$ws.Range("G14").Value = 75.625
$ws.Range("G20").Value = 75.625

# --- 2. Apply Text number format to columns A, B, G, H, I ------------------
$ws.Range("A1:A25").NumberFormat = "@"
$ws.Range("B1:B25").NumberFormat = "@"
$ws.Range("G1:H25").NumberFormat = "@"
$ws.Range("I1").NumberFormat = "@"

# --- 3. Column widths for the (previously implicit) default-width columns --
# Columns B and H now carry explicit column formatting; re-assert their
# (unchanged) default width so the column metadata is written out.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth()
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(8).ColumnWidth()

# --- 4. Cursor / selection position ----------------------------------------
$ws.Range("J18").Select()
